$d = $word.ActiveDocument
$replacements = @(
    @{ Old = "[[PERSON_30]] – „k [[PERSON_31]]“, „o [[PERSON_31]]“"; New = "[[PERSON_30]] – „k [[PERSON_30]]“, „o [[PERSON_30]]“" }
    @{ Old = "[[PERSON_32]] – „o [[PERSON_32]]“, „s [[PERSON_32]]“"; New = "[[PERSON_31]] – „o [[PERSON_31]]“, „s [[PERSON_31]]“" }
    @{ Old = "[[PERSON_33]] – „s [[PERSON_33]]“, „o [[PERSON_33]]“"; New = "[[PERSON_32]] – „s [[PERSON_32]]“, „o [[PERSON_32]]“" }
    @{ Old = "[[PERSON_34]] – „s [[PERSON_34]]“, „o [[PERSON_34]]“"; New = "[[PERSON_33]] – „s [[PERSON_33]]“, „o [[PERSON_33]]“" }
    @{ Old = "[[PERSON_35]] – „k [[PERSON_35]]“, „s [[PERSON_35]]“"; New = "[[PERSON_34]] – „k [[PERSON_34]]“, „s [[PERSON_34]]“" }
    @{ Old = "[[PERSON_36]] – „pro [[PERSON_36]]“, „o [[PERSON_36]]“"; New = "[[PERSON_35]] – „pro [[PERSON_35]]“, „o [[PERSON_35]]“" }
    @{ Old = "[[PERSON_37]] – „k [[PERSON_37]]“, „o [[PERSON_37]]“"; New = "[[PERSON_36]] – „k [[PERSON_36]]“, „o [[PERSON_36]]“" }
    @{ Old = "[[PERSON_38]] – „o [[PERSON_38]]“, „s [[PERSON_38]]“"; New = "[[PERSON_37]] – „o [[PERSON_37]]“, „s [[PERSON_37]]“" }
    @{ Old = "[[PERSON_39]] – „s [[PERSON_39]]“, „o [[PERSON_39]]“"; New = "[[PERSON_38]] – „s [[PERSON_38]]“, „o [[PERSON_38]]“" }
    @{ Old = "[[PERSON_40]] – „s [[PERSON_40]]“, „o [[PERSON_40]]“"; New = "[[PERSON_39]] – „s [[PERSON_39]]“, „o [[PERSON_39]]“" }
    @{ Old = "[[PERSON_41]] – „u [[PERSON_42]]“, „o [[PERSON_41]]“"; New = "[[PERSON_40]] – „u [[PERSON_40]]“, „o [[PERSON_40]]“" }
    @{ Old = "[[PERSON_43]] – „se [[PERSON_43]]“, „o [[PERSON_43]]“"; New = "[[PERSON_41]] – „se [[PERSON_41]]“, „o [[PERSON_41]]“" }
    @{ Old = "[[PERSON_44]] – „o [[PERSON_44]]“, „s [[PERSON_44]]“"; New = "[[PERSON_42]] – „o [[PERSON_42]]“, „s [[PERSON_42]]“" }
    @{ Old = "[[PERSON_45]] – „k [[PERSON_45]]“, „o [[PERSON_45]]“"; New = "[[PERSON_43]] – „k [[PERSON_43]]“, „o [[PERSON_43]]“" }
    @{ Old = "[[PERSON_46]] – „o [[PERSON_47]]“, „s [[PERSON_46]]“"; New = "[[PERSON_44]] – „o [[PERSON_44]]“, „s [[PERSON_44]]“" }
    @{ Old = "[[PERSON_48]] – „s [[PERSON_48]]“, „o [[PERSON_48]]“"; New = "[[PERSON_45]] – „s [[PERSON_45]]“, „o [[PERSON_45]]“" }
    @{ Old = "[[PERSON_49]] – „s [[PERSON_49]]“, „o [[PERSON_49]]“"; New = "[[PERSON_46]] – „s [[PERSON_46]]“, „o [[PERSON_46]]“" }
    @{ Old = "[[PERSON_50]] – „o [[PERSON_50]]“, „s [[PERSON_50]]“"; New = "[[PERSON_47]] – „o [[PERSON_47]]“, „s [[PERSON_47]]“" }
    @{ Old = "[[PERSON_51]] – „s [[PERSON_51]]“, „o [[PERSON_51]]“"; New = "[[PERSON_48]] – „s [[PERSON_48]]“, „o [[PERSON_48]]“" }
    @{ Old = "[[PERSON_52]] – „o [[PERSON_52]]“, „s [[PERSON_52]]“"; New = "[[PERSON_49]] – „o [[PERSON_49]]“, „s [[PERSON_49]]“" }
    @{ Old = "[[PERSON_53]] – „s [[PERSON_53]]“, „o [[PERSON_53]]“"; New = "[[PERSON_50]] – „s [[PERSON_50]]“, „o [[PERSON_50]]“" }
)

foreach ($rep in $replacements) {
    $range = $d.Content
    $ok = $range.Find.Execute($rep.Old, $true, $true, $false, $false, $false, $true, 1, $false, $rep.New, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $($rep.Old)"
    }
}
Write-Output "Done"